$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1455880165028277
$ws.Range("C2").Value = 0.1280610774805511
$ws.Range("D2").Value = 0.09270042616372805
$ws.Range("B3").Value = 0.09139151174382656
$ws.Range("C3").Value = 2.548199582842637
$ws.Range("D3").Value = 0.1184211724430862
$ws.Range("B4").Value = 0.05983318338467544
$ws.Range("C4").Value = 4.149682042742901
$ws.Range("D4").Value = 0.1067871507012871
$ws.Range("B5").Value = 0.108364248077063
$ws.Range("C5").Value = 4.935671882245719
$ws.Range("D5").Value = 0.1606575392949604
$ws.Range("B6").Value = 0.06804098872428861
$ws.Range("C6").Value = 4.707338191198275
$ws.Range("D6").Value = 0.201327409164833
$ws.Range("B7").Value = 0.09206696030400119
$ws.Range("C7").Value = 4.201761483090537
$ws.Range("D7").Value = 0.1025373591177245
$ws.Range("B8").Value = 0.09046146085689701
$ws.Range("C8").Value = 2.778719414675342
$ws.Range("D8").Value = 0.1558098794300889
$ws.Range("B9").Value = 0.07886436601773093
$ws.Range("C9").Value = 0.03350573213367947
$ws.Range("D9").Value = 0.154601706930968
$ws.Range("B10").Value = 0.09506032031438941
$ws.Range("C10").Value = -1.495458666213109
$ws.Range("D10").Value = 0.1650805057731569
$ws.Range("B11").Value = 0.07869772232140226
$ws.Range("C11").Value = -3.653063522703919
$ws.Range("D11").Value = 0.2037632037382193
$ws.Range("B12").Value = 0.103082925621828
$ws.Range("C12").Value = -4.46875886981198
$ws.Range("D12").Value = 0.1469853422132432
$ws.Range("B13").Value = 0.1388981845366059
$ws.Range("C13").Value = -4.858232231565448
$ws.Range("D13").Value = 0.1858437507872601
$ws.Range("B14").Value = 0.09069162671203221
$ws.Range("C14").Value = -4.450859129078268
$ws.Range("D14").Value = 0.213201029001834
$ws.Range("B15").Value = 0.1115689197831609
$ws.Range("C15").Value = -2.986969868350579
$ws.Range("D15").Value = 0.1278059581188271
$ws.Range("B16").Value = 0.09245792824833171
$ws.Range("C16").Value = -1.547766147877984
$ws.Range("D16").Value = 0.1629921861483814
$ws.Range("B17").Value = 0.1425707420760201
$ws.Range("C17").Value = 1.024787522364725
$ws.Range("D17").Value = 0.09930502171121953
$ws.Range("B18").Value = 0.09094871485179887
$ws.Range("C18").Value = 2.958544685266378
$ws.Range("D18").Value = 0.1409829891313443
$ws.Range("B19").Value = 0.1284462752495574
$ws.Range("C19").Value = 4.375785770475057
$ws.Range("D19").Value = 0.135631851107491
$ws.Range("B20").Value = 0.1066601377215915
$ws.Range("C20").Value = 5.05021771817742
$ws.Range("D20").Value = 0.195627997566301
$ws.Range("B21").Value = 0.05360559235318207
$ws.Range("C21").Value = 4.611108772339825
$ws.Range("D21").Value = 0.09412752215797995
$ws.Range("B22").Value = 0.1330634175891486
$ws.Range("C22").Value = 3.489665551903669
$ws.Range("D22").Value = 0.1275693512478061
$ws.Range("B23").Value = 0.1470496085731526
$ws.Range("C23").Value = 1.61533540159829
$ws.Range("D23").Value = 0.11133842161863
$ws.Range("B24").Value = 0.08379739931321023
$ws.Range("C24").Value = -0.04494584946378549
$ws.Range("D24").Value = 0.1529469426678755
$ws.Range("B25").Value = 0.06029211290147602
$ws.Range("C25").Value = -1.942407820099582
$ws.Range("D25").Value = 0.1616473544859163
$ws.Range("B26").Value = 0.0880550385632594
$ws.Range("C26").Value = -3.655319961100368
$ws.Range("D26").Value = 0.2012530015560313
$ws.Range("B27").Value = 0.0652092905486486
$ws.Range("C27").Value = -5.004316611494589
$ws.Range("D27").Value = 0.1945982067478729
$ws.Range("B28").Value = 0.05023201584900652
$ws.Range("C28").Value = -5.082250238617195
$ws.Range("D28").Value = 0.07613404961077201
$ws.Range("B29").Value = 0.07485939915551286
$ws.Range("C29").Value = -3.93128332438776
$ws.Range("D29").Value = 0.1383974990886414
$ws.Range("B30").Value = 0.1288521967776751
$ws.Range("C30").Value = -2.169156619365828
$ws.Range("D30").Value = 0.1204428614699008
$ws.Range("B31").Value = 0.05598252222373116
$ws.Range("C31").Value = 0.1443862475804094
$ws.Range("D31").Value = 0.1934376238137063
